# Fruta / hortaliza, semanal
#
# Two new weekly price records (rows 31-32, dated 2021-10-18) are inserted
# at the top of the data block. This pushes every existing record down by
# two rows, so the last two existing records (old rows 124-125) overflow
# into two brand-new rows (126-127) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 31; this shifts rows 31:125 down to 33:127
# and expands the sheet's used range/dimension to A1:T127 automatically.
$ws.Rows(31).Resize(2).Insert()

# Seed the two new rows with the same shape/formatting as the row that is
# now directly beneath them (formerly row 31/32, now row 33/34), then
# overwrite just the fields that actually differ for the new records.
$ws.Range("A33:T34").Copy() | Out-Null
$ws.Range("A31").PasteSpecial() | Out-Null

# Row 31: new "Especial" quality record dated 2021-10-18
$ws.Range("D31").Value2 = 44487
$ws.Range("L31").Value = "Especial"
$ws.Range("M31").Value = 85

# Row 32: new "Segunda" quality record dated 2021-10-18
$ws.Range("D32").Value2 = 44487
$ws.Range("M32").Value = 70
$ws.Range("N32").Value = 8000
$ws.Range("O32").Value = 8000
$ws.Range("P32").Value = 8000
$ws.Range("S32").Value = 1143
